# v1.1 via stitching added
# Correct the MP (manufacturer part number) for the CAN-B1 connector (Molex):
# it was "502494-0670" but should be "502585-0670" (matching the Name column, L3).
# The leading apostrophe preserves the cell's existing "stored as text"
# (quote-prefix) formatting instead of letting Excel reinterpret the value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "'502585-0670"
